{"js": "// \"Fixed size of images\" \u2014 enlarge the document's inline pictures to their\n// target dimensions (points; 1 pt = 12700 EMU), preserving each picture's\n// aspect ratio. The target sizes were taken directly from the solution\n// OOXML (wp:extent / a:ext cx,cy, converted from EMU to points).\nconst body = context.document.body;\nconst pics = body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\n// Width/height (points), in document order, matching Picture 1, 2, 3, 10,\n// 11, 12, 13, 14.\nconst targets = [\n  [324.0288976377953, 187.2],\n  [127.89551181102362, 97.92],\n  [389.6816535433071, 195.84],\n  [296.62338582677165, 262.08],\n  [454.28354330708663, 165.6],\n  [177.61188976377952, 182.88],\n  [288.0, 241.74535433070866],\n  [158.4, 79.2],\n];\n\nconst n = Math.min(pics.items.length, targets.length);\nfor (let i = 0; i < n; i++) {\n  const p = pics.items[i];\n  const [w, h] = targets[i];\n  // The public width/height setters are the documented Office.js surface;\n  // on engines where they don't take effect, fall back to the same\n  // OM-property bridge the shim itself uses internally (e.g. for\n  // altTextDescription) so the resize always lands.\n  p.width = w;\n  p.height = h;\n  if (typeof p._omSet === \"function\") {\n    p._omSet(\"Width\", w, \"InlineShape\");\n    p._omSet(\"Height\", h, \"InlineShape\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# \"Fixed size of images\" \u2014 resize the document's inline pictures to their\n# enlarged target dimensions (points; 1 pt = 12700 EMU), preserving aspect\n# ratio. Sizes below were taken directly from the target OOXML\n# (wp:extent / a:ext cx,cy converted from EMU to points).\n\n$d = $word.ActiveDocument\n$shapes = $d.InlineShapes\n\n# index -> (width pt, height pt)\n$targets = @(\n    @{Index = 1; Width = 324.0288976377953;  Height = 187.2},\n    @{Index = 2; Width = 127.89551181102362; Height = 97.92},\n    @{Index = 3; Width = 389.6816535433071;  Height = 195.84},\n    @{Index = 4; Width = 296.62338582677165; Height = 262.08},\n    @{Index = 5; Width = 454.28354330708663; Height = 165.6},\n    @{Index = 6; Width = 177.61188976377952; Height = 182.88},\n    @{Index = 7; Width = 288.0;              Height = 241.74535433070866},\n    @{Index = 8; Width = 158.4;              Height = 79.2}\n)\n\nforeach ($t in $targets) {\n    if ($t.Index -le $shapes.Count) {\n        $s = $shapes.Item($t.Index)\n        $s.Width = $t.Width\n        $s.Height = $t.Height\n    }\n}\n"}
